$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8356053233146667
$ws.Range("B1").Value = 0.9433885812759399
$ws.Range("C1").Value = 0.7396496534347534
$ws.Range("D1").Value = 0.6909381151199341
$ws.Range("E1").Value = 0.7210728526115417
